# Update cryptocurrency price/volume data to reflect the latest scrape.
# Mirrors the GitHub Actions "Updated cryptos list" commit.
#
# Price/volume values are written as plain text (matching the workbook's
# existing inline-string convention) even when they look numeric, so cells
# such as "211.46" or "0.247" are not silently reinterpreted as floating
# point numbers by Excel's auto-detection. NumberFormat is forced to Text
# ("@") just long enough to commit the literal string, then ClearFormats()
# removes the temporary formatting so the cell keeps its original
# (default/general) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.659.75"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "1.597.86"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.46"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.247"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.59"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "1.822.03"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.03"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.552.81"
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.03"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "26.637.62"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "0.0₃0738"
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.92"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.03"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.47%  "
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.99"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.34"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.27"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0517"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.40%  "
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("D34").Value = "1.287.98"
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.614"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -8.58%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("E40").Value = "  +18.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.50"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.08%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.51"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").Value = "1.734.24"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.76"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.79%  "
$ws.Range("E47").Value = "  -3.37%  "
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0509"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.37"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.10%  "
